# modified test cases on overdue fix
#
# Summary sheet: the expected "over due" amounts move from 117.19 to
# 117.25 (A3/E3), F3 keeps its value but all three cells lose their
# custom "0.00" number format (style 22) in favour of the sheet's plain
# "General" style (style 14).
#
# Repayment schedule sheet: the schedule is recalculated so that a 7th
# instalment period is now required (a new row 9), the 6th period's
# figures change (late payment no longer clears the loan), and almost
# every numeric cell's bespoke "0.00"/"#,##0.00" number format (styles
# 22/23) is replaced by the sheet's plain style (14), except for the
# running "Loan Balance" column G, which keeps a numeric format (style
# 16 = "#,##0" or 17 = "#,##0.00" depending on the row). The now-unused
# trailing column R is cleared out entirely.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Summary")

$ws2.Range("A3").Value2 = 117.25
$ws2.Range("E3").Value2 = 117.25

# Re-format A3/E3/F3 from style 22 ("0.00") to the sheet's plain style 14
# (General) by copying the format off an already-correct neighbour cell.
$ws2.Range("B3").Copy()
$ws2.Range("A3,E3,F3").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws2.Activate()
$ws2.Range("A8:XFD15").Select()

# ---------------------------------------------------------------------
# Repayment schedule sheet
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Repayment schedule")

# Drop the stray column R entirely (content + formatting).
$ws3.Range("R3:R8").Clear()

# Normalise every style-22/23 cell (plus the never-styled A/B/I/J cells on
# rows 4-8) onto the sheet's plain style 14, by copying format from an
# already-style-14 cell (A2).
$ws3.Range("A2").Copy()
$ws3.Range("F2,H2,P2,F3,H3,K3,P3,A4:B4,D4:E4,F4,H4,I4:J4,K4,P4,A5:B5,D5:E5,F5,H5,I5:J5,K5,P5,A6:B6,D6:E6,F6,H6,I6:J6,K6,P6,A7:B7,D7:E7,F7,G7,H7,I7:J7,K7,P7,A8:B8,D8:E8,F8,G8,H8,I8:J8,K8,P8").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Loan-balance column G keeps a numeric format: style 16 ("#,##0") for
# G2/G6, style 17 ("#,##0.00") for G3/G4/G5.
$ws2.Range("A2").Copy()
$ws3.Range("G2,G6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws2.Range("F2").Copy()
$ws3.Range("G3:G5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Value updates for the recalculated schedule.
$ws3.Range("P2").ClearContents()

$ws3.Range("G6").Value2 = 1690

$ws3.Range("G7").Value2 = 852.58

$ws3.Range("F8").Value2 = 844.4
$ws3.Range("G8").Value2 = 8.18
$ws3.Range("H8").Value2 = 7.1
$ws3.Range("K8").Value2 = 851.5
$ws3.Range("P8").Value2 = 851.5

# New 7th instalment row.
$ws3.Range("A2:P2").Copy()
$ws3.Range("A9:P9").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws3.Range("A9").Value2 = 7
$ws3.Range("B9").Value2 = 31
$ws3.Range("C9").Value2 = 42217
$ws3.Range("F9").Value2 = 8.18
$ws3.Range("G9").Value2 = 0
$ws3.Range("H9").Value2 = 0.07
$ws3.Range("I9").Value2 = 0
$ws3.Range("J9").Value2 = 0
$ws3.Range("K9").Value2 = 8.25
$ws3.Range("L9").Value2 = 0
$ws3.Range("M9").Value2 = 0
$ws3.Range("N9").Value2 = 0
$ws3.Range("O9").Value2 = 0
$ws3.Range("P9").Value2 = 8.25

$ws3.Activate()
$ws3.Range("A10:XFD10").Select()
